$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for columns D, J, K, L, M, P
$ws.Range("D2").Value2 = 44749
$ws.Range("J2").Value2 = 90
$ws.Range("K2").Value2 = 17000
$ws.Range("L2").Value2 = 18000
$ws.Range("M2").Value2 = 17556
$ws.Range("P2").Value2 = 1170

$ws.Range("D3").Value2 = 44750
$ws.Range("J3").Value2 = 140
$ws.Range("K3").Value2 = 19000
$ws.Range("L3").Value2 = 20000
$ws.Range("M3").Value2 = 19571
$ws.Range("P3").Value2 = 1305
